$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(177, 1).Value = 176
$ws.Cells.Item(177, 2).Value = "Pós OS"
$ws.Cells.Item(177, 3).Value = 8005258758
$ws.Cells.Item(177, 4).Value = 46023
$ws.Cells.Item(177, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(177, 6).Value = "Detratores"
$ws.Cells.Item(177, 7).Value = "Tecnico chegou num horário diferente do agendado. Eu estava de saida para outro compromisso e me atrasou. Sugiro que os agendamentos sejam cumpridos conforme o cliente demanda."
$ws.Cells.Item(177, 8).Value = "Campo"
$ws.Cells.Item(177, 9).Value = "Fora do período agendado"

$ws.Cells.Item(178, 1).Value = 177
$ws.Cells.Item(178, 2).Value = "Pós OS"
$ws.Cells.Item(178, 3).Value = 8005224647
$ws.Cells.Item(178, 4).Value = 46023
$ws.Cells.Item(178, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(178, 6).Value = "Detratores"
$ws.Cells.Item(178, 7).Value = "O técnico não compareceu. Não houve visita. Essa é a terceira vez seguida que acontece isso. Extremamente insatisfeito"
$ws.Cells.Item(178, 8).Value = "Campo"
$ws.Cells.Item(178, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(179, 1).Value = 178
$ws.Cells.Item(179, 2).Value = "Instalação"
$ws.Cells.Item(179, 3).Value = 8005267106
$ws.Cells.Item(179, 4).Value = 46023
$ws.Cells.Item(179, 5).Value = "FRQ_ECO_SP_GUARULHOS_2"
$ws.Cells.Item(179, 6).Value = "Detratores"
$ws.Cells.Item(179, 7).Value = "Muito ruim. O instalador chegou de manhã mesmo combinando a tarde. Tocou o interfone e acordou minha mae de 85 anos. E deixou um serviço porco. Fio aparecendo.  Péssimo.  Uma vergonha"
$ws.Cells.Item(179, 8).Value = "Campo"
$ws.Cells.Item(179, 9).Value = "Fora do período agendado"

$ws.Cells.Item(180, 1).Value = 179
$ws.Cells.Item(180, 2).Value = "Pós OS"
$ws.Cells.Item(180, 3).Value = 8005234424
$ws.Cells.Item(180, 4).Value = 46023
$ws.Cells.Item(180, 5).Value = "FRQ_ECO_SP_S B CAMPO"
$ws.Cells.Item(180, 6).Value = "Detratores"
$ws.Cells.Item(180, 7).Value = "Houve demora para manutenção. O fluxo de água do meu aparelho é baixo."
$ws.Cells.Item(180, 8).Value = "Qualidade do Produto"
$ws.Cells.Item(180, 9).Value = "Vazão/pressão da água"

$ws.Cells.Item(181, 1).Value = 180
$ws.Cells.Item(181, 2).Value = "Pós OS"
$ws.Cells.Item(181, 3).Value = 8005271095
$ws.Cells.Item(181, 4).Value = 46023
$ws.Cells.Item(181, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(181, 6).Value = "Detratores"
$ws.Cells.Item(181, 7).Value = "Pela quarta vez foi agendada a manutenção e o técnico não veio. Não há como avaliar algo que não aconteceu Inclusive solicito que verifiquem minha situação, pois é quarto agendamento que o técnico não vem e não recebi nenhuma explicação"
$ws.Cells.Item(181, 8).Value = "Campo"
$ws.Cells.Item(181, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(182, 1).Value = 181
$ws.Cells.Item(182, 2).Value = "Pós OS"
$ws.Cells.Item(182, 3).Value = 8005260391
$ws.Cells.Item(182, 4).Value = 46023
$ws.Cells.Item(182, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(182, 6).Value = "Detratores"
$ws.Cells.Item(182, 7).Value = "Já fiz duas chamadas de técnico e não apareceu! Fico esperando, perdendo trabalho e A Brastemp não apareceu"
$ws.Cells.Item(182, 8).Value = "Campo"
$ws.Cells.Item(182, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(183, 1).Value = 182
$ws.Cells.Item(183, 2).Value = "Pós OS"
$ws.Cells.Item(183, 3).Value = 8005218070
$ws.Cells.Item(183, 4).Value = 46023
$ws.Cells.Item(183, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(183, 6).Value = "Detratores"
$ws.Cells.Item(183, 7).Value = "Mais uma vez o técnico não veio . Estou exausta de reclamar. Já liguei, tentei cancelar . A atendente conseguiu baixar minha parcela, porém ainda com problemas em vir aqui e trocar o filtro. Já faltei trabalho, já deixei de ir a tratamentos médicos . Mas eu tenho compromisso. Nem satisfação me dão. Ou…. Chegam aqui em horários que eu não estou. Simplesmente assim.  Se eu deixasse de pagar, provavelmente muito rapidamente me cobrariam. Mas como eu pago em dia me tratam assim. Sinceramente ."
$ws.Cells.Item(183, 8).Value = "Campo"
$ws.Cells.Item(183, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(184, 1).Value = 183
$ws.Cells.Item(184, 2).Value = "Pós OS"
$ws.Cells.Item(184, 3).Value = 8005267385
$ws.Cells.Item(184, 4).Value = 46023
$ws.Cells.Item(184, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(184, 6).Value = "Detratores"
$ws.Cells.Item(184, 7).Value = "Fiz a reclamação de nao terem feito a instalação do gas para agua, ficaram de vir no sabado no horário da manhã e de novo ninguem apareceu. Estou achando pessimo pq quis locar por conta da agua com gas e até agora nao tenho."
$ws.Cells.Item(184, 8).Value = "Campo"
$ws.Cells.Item(184, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(185, 1).Value = 184
$ws.Cells.Item(185, 2).Value = "Pós OS"
$ws.Cells.Item(185, 3).Value = 8005247512
$ws.Cells.Item(185, 4).Value = 46023
$ws.Cells.Item(185, 5).Value = "FRQ_ECO_SP_ZONA_SUL_03"
$ws.Cells.Item(185, 6).Value = "Neutros"
$ws.Cells.Item(185, 7).Value = "Caro, só reduz o preço depois da minha reclamação.Precisa mandar um lembrete um dia antes da manutençãoA forma de agendamento é irritante. Podia ser um whats. É horrível receber q ligação"
$ws.Cells.Item(185, 8).Value = "Outros"
$ws.Cells.Item(185, 9).Value = "Preço elevado"

$ws.Cells.Item(186, 1).Value = 185
$ws.Cells.Item(186, 2).Value = "Pós OS"
$ws.Cells.Item(186, 3).Value = 8005190250
$ws.Cells.Item(186, 4).Value = 46023
$ws.Cells.Item(186, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(186, 6).Value = "Detratores"
$ws.Cells.Item(186, 7).Value = "_ Estou pedindo a visita do técnico para trocar o filtro do meu purificador e ele nunca vem. Tenho vários protocolos de pedido, e ninguém aparece:05/11/25 - Na parte da Tarde - entre 13:00 e 18:00hs - Ninguém apareceu Ordem de Serviço - N° 800519025016/01/26 - Na parte da Tarde - entre 13:00 e 18:00hs - Ninguém apareceu Ordem de Serviço - N° 8005190250Falta de profissionalismo, pessoal descompromissadoPago em dia as mensalidades do Purificador BRASTEMP, porém a Empresa não consegue fazer as manutenções no período combinado."
$ws.Cells.Item(186, 8).Value = "Campo"
$ws.Cells.Item(186, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(187, 1).Value = 186
$ws.Cells.Item(187, 2).Value = "Pós OS"
$ws.Cells.Item(187, 3).Value = 8005267277
$ws.Cells.Item(187, 4).Value = 46023
$ws.Cells.Item(187, 5).Value = "FRQ_ECO_SP_OSASCO"
$ws.Cells.Item(187, 6).Value = "Detratores"
$ws.Cells.Item(187, 7).Value = "O serviço do suporte não atendeu minhas expectativas"
$ws.Cells.Item(187, 8).Value = "Outros"
$ws.Cells.Item(187, 9).Value = "Insatisfação geral"

$ws.Cells.Item(188, 1).Value = 187
$ws.Cells.Item(188, 2).Value = "Pós OS"
$ws.Cells.Item(188, 3).Value = 8005262536
$ws.Cells.Item(188, 4).Value = 46023
$ws.Cells.Item(188, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(188, 6).Value = "Detratores"
$ws.Cells.Item(188, 7).Value = "Mais de 20 dias sem água gelada e a Brastemp não está  nem aí com o consumidor, mesmo ele sendo assinante mais de 10 anos. Uma vergonha."
$ws.Cells.Item(188, 8).Value = "Capacidade"
$ws.Cells.Item(188, 9).Value = "Data Distante"

$ws.Cells.Item(189, 1).Value = 188
$ws.Cells.Item(189, 2).Value = "Pós OS"
$ws.Cells.Item(189, 3).Value = 8005263939
$ws.Cells.Item(189, 4).Value = 46023
$ws.Cells.Item(189, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(189, 6).Value = "Detratores"
$ws.Cells.Item(189, 7).Value = "Técnico não veio e não deu satisfação."
$ws.Cells.Item(189, 8).Value = "Campo"
$ws.Cells.Item(189, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(190, 1).Value = 189
$ws.Cells.Item(190, 2).Value = "Pós OS"
$ws.Cells.Item(190, 3).Value = 8005238421
$ws.Cells.Item(190, 4).Value = 46023
$ws.Cells.Item(190, 5).Value = "AT_ECO_CE_FORTALEZA"
$ws.Cells.Item(190, 6).Value = "Detratores"
$ws.Cells.Item(190, 7).Value = "Não estou gostando, o problema sempre volta.O técnico saiu e o problema persistiu e ficou dando choque. Agora é que não quero mais.Quero entregar o aparelho e findar o meu aluguel."
$ws.Cells.Item(190, 8).Value = "Campo"
$ws.Cells.Item(190, 9).Value = "Reincidência"

$ws.Cells.Item(191, 1).Value = 190
$ws.Cells.Item(191, 2).Value = "Pós OS"
$ws.Cells.Item(191, 3).Value = 8005266918
$ws.Cells.Item(191, 4).Value = 46023
$ws.Cells.Item(191, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(191, 6).Value = "Detratores"
$ws.Cells.Item(191, 7).Value = "Mais uma vez.. agendaram e não vieram"
$ws.Cells.Item(191, 8).Value = "Campo"
$ws.Cells.Item(191, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(192, 1).Value = 191
$ws.Cells.Item(192, 2).Value = "Pós OS"
$ws.Cells.Item(192, 3).Value = 8005273027
$ws.Cells.Item(192, 4).Value = 46023
$ws.Cells.Item(192, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(192, 6).Value = "Detratores"
$ws.Cells.Item(192, 7).Value = "O problema não foi resolvido."
$ws.Cells.Item(192, 8).Value = "Campo"
$ws.Cells.Item(192, 9).Value = "Reincidência"

$ws.Cells.Item(193, 1).Value = 192
$ws.Cells.Item(193, 2).Value = "Pós OS"
$ws.Cells.Item(193, 3).Value = 8005251083
$ws.Cells.Item(193, 4).Value = 46023
$ws.Cells.Item(193, 5).Value = "FRQ_ECO_SP_ZONA_SUL_03"
$ws.Cells.Item(193, 6).Value = "Neutros"
$ws.Cells.Item(193, 7).Value = "A disponibilidade de uma visita próxima foi ruim, passei quase 3 semanas com O purificador com cheiro e gosto ruim, porque não tinham disponibilidade antes."
$ws.Cells.Item(193, 8).Value = "Capacidade"
$ws.Cells.Item(193, 9).Value = "Data Distante"

$ws.Cells.Item(194, 1).Value = 193
$ws.Cells.Item(194, 2).Value = "Pós OS"
$ws.Cells.Item(194, 3).Value = 8005248857
$ws.Cells.Item(194, 4).Value = 46023
$ws.Cells.Item(194, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(194, 6).Value = "Detratores"
$ws.Cells.Item(194, 7).Value = "Não houve visita tecnica"
$ws.Cells.Item(194, 8).Value = "Campo"
$ws.Cells.Item(194, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(195, 1).Value = 194
$ws.Cells.Item(195, 2).Value = "Pós OS"
$ws.Cells.Item(195, 3).Value = 8005207126
$ws.Cells.Item(195, 4).Value = 46023
$ws.Cells.Item(195, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(195, 6).Value = "Detratores"
$ws.Cells.Item(195, 7).Value = "Não tive visita , o técnico não veio , estou muito insatisfeita , segunda vez agendo e nada"
$ws.Cells.Item(195, 8).Value = "Campo"
$ws.Cells.Item(195, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(196, 1).Value = 195
$ws.Cells.Item(196, 2).Value = "Pós OS"
$ws.Cells.Item(196, 3).Value = 8005267788
$ws.Cells.Item(196, 4).Value = 46023
$ws.Cells.Item(196, 5).Value = "FRQ_ECO_SP_OSASCO"
$ws.Cells.Item(196, 6).Value = "Neutros"
$ws.Cells.Item(196, 7).Value = "Na minha opinião o filtro deverá ser trocada semestralmente, principalmente pela qualidade d'água Sabesp tem piorado e muito, isso naturalmente exige mais do filtro e a cada 6 meses o técnico vem e diz que a agua esta bom, e quando peço para trocar eles não carregam esse filtro adicional, tenho que solicitar para dizer que a agua não esta boa mesmo. Acabo comprando agua no supermercado, porque além do cheiro desagradável apresenta coloração."
$ws.Cells.Item(196, 8).Value = "Campo"
$ws.Cells.Item(196, 9).Value = "Qualidade da manutenção"

$ws.Cells.Item(197, 1).Value = 196
$ws.Cells.Item(197, 2).Value = "Instalação"
$ws.Cells.Item(197, 3).Value = 8005262682
$ws.Cells.Item(197, 4).Value = 46023
$ws.Cells.Item(197, 5).Value = "FRQ_ECO_SP_S B CAMPO"
$ws.Cells.Item(197, 6).Value = "Detratores"
$ws.Cells.Item(197, 7).Value = "Não funcionou na instalação, solicitei novamente o técnico e o mesmo não veio."
$ws.Cells.Item(197, 8).Value = "Campo"
$ws.Cells.Item(197, 9).Value = "Qualidade da instalação"

$ws.Cells.Item(198, 1).Value = 197
$ws.Cells.Item(198, 2).Value = "Pós OS"
$ws.Cells.Item(198, 3).Value = 8005261464
$ws.Cells.Item(198, 4).Value = 46023
$ws.Cells.Item(198, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(198, 6).Value = "Detratores"
$ws.Cells.Item(198, 7).Value = "O técnico não compareceu!"
$ws.Cells.Item(198, 8).Value = "Campo"
$ws.Cells.Item(198, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(199, 1).Value = 198
$ws.Cells.Item(199, 2).Value = "Pós OS"
$ws.Cells.Item(199, 3).Value = 8005250369
$ws.Cells.Item(199, 4).Value = 46023
$ws.Cells.Item(199, 5).Value = "FRQ_ECO_SP_ZONA_SUL_03"
$ws.Cells.Item(199, 6).Value = "Detratores"
$ws.Cells.Item(199, 7).Value = "O técnico sugeriu, pela segunda vez, não trocar o filtro. Devido a nossa insistência, uma vez que na visita anterior a troca não foi feita, ele realizou a troca. Ou seja, qual é o critério técnico desta visita técnica? Talvez eu seja um dos clientes mais antigos do Purificador Brastemp, pela primeira vez, estou em dúvida sobre a qualidade do serviço."
$ws.Cells.Item(199, 8).Value = "Campo"
$ws.Cells.Item(199, 9).Value = "Qualidade da manutenção"

$ws.Cells.Item(200, 1).Value = 199
$ws.Cells.Item(200, 2).Value = "Pós OS"
$ws.Cells.Item(200, 3).Value = 8005243013
$ws.Cells.Item(200, 4).Value = 46023
$ws.Cells.Item(200, 5).Value = "FRQ_ECO_SP_S B CAMPO"
$ws.Cells.Item(200, 6).Value = "Detratores"
$ws.Cells.Item(200, 7).Value = "Não foi feito o combinado em instalar o MIB geo"
$ws.Cells.Item(200, 8).Value = "Campo"
$ws.Cells.Item(200, 9).Value = "Qualidade da manutenção"

$ws.Cells.Item(201, 1).Value = 200
$ws.Cells.Item(201, 2).Value = "Pós OS"
$ws.Cells.Item(201, 3).Value = 8005250128
$ws.Cells.Item(201, 4).Value = 46023
$ws.Cells.Item(201, 5).Value = "FRQ_ECO_SP_ZONA_SUL_03"
$ws.Cells.Item(201, 6).Value = "Detratores"
$ws.Cells.Item(201, 7).Value = "Muito simples a limpeza, nem ao menos um produtos tipo para desinfetar o filtro todo após o manuseio."
$ws.Cells.Item(201, 8).Value = "Campo"
$ws.Cells.Item(201, 9).Value = "Qualidade da manutenção"

$ws.Cells.Item(202, 1).Value = 201
$ws.Cells.Item(202, 2).Value = "Pós OS"
$ws.Cells.Item(202, 3).Value = 8005257348
$ws.Cells.Item(202, 4).Value = 46023
$ws.Cells.Item(202, 5).Value = "FRQ_ECO_RJ_OESTE"
$ws.Cells.Item(202, 6).Value = "Detratores"
$ws.Cells.Item(202, 7).Value = "Não houve a visita."
$ws.Cells.Item(202, 8).Value = "Campo"
$ws.Cells.Item(202, 9).Value = "Técnico não cumpriu a agenda"

$ws.Cells.Item(203, 1).Value = 202
$ws.Cells.Item(203, 2).Value = "Pós OS"
$ws.Cells.Item(203, 3).Value = 8005250948
$ws.Cells.Item(203, 4).Value = 46023
$ws.Cells.Item(203, 5).Value = "FRQ_ECO_RS_PORTO ALEGRE_2"
$ws.Cells.Item(203, 6).Value = "Detratores"
$ws.Cells.Item(203, 7).Value = "depois de 5 agendamentos e muitas ligações ocorreu o serviço"
$ws.Cells.Item(203, 8).Value = "Capacidade"
$ws.Cells.Item(203, 9).Value = "Data Distante"

$ws.Range("A177:I203").Select()
